$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.622.60'
$ws.Range("E2").Value = '  +0.83%  '
$ws.Range("D3").Value = '1.801.20'
$ws.Range("E3").Value = '  -0.99%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.06'
$ws.Range("E5").Value = '  -1.46%  '
$ws.Range("E6").Value = '  -3.12%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '33.04'
$ws.Range("E8").Value = '  +4.60%  '
$ws.Range("E9").Value = '  -0.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0667'
$ws.Range("E10").Value = '  -1.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0933'
$ws.Range("E11").Value = '  +0.11%  '
$ws.Range("D12").Value = '2.054.83'
$ws.Range("E12").Value = '  -1.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.12'
$ws.Range("E13").Value = '  +8.50%  '
$ws.Range("D14").Value = '1.790.11'
$ws.Range("E14").Value = '  -1.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.644'
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("D16").Value = '34.546.58'
$ws.Range("E16").Value = '  +0.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.29'
$ws.Range("E17").Value = '  +0.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.74'
$ws.Range("E18").Value = '  -0.98%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '256.29'
$ws.Range("E19").Value = '  -1.02%  '
$ws.Range("D20").Value = '0.0₃0760'
$ws.Range("E20").Value = '  +0.59%  '
$ws.Range("E21").Value = '  -0.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.49'
$ws.Range("E22").Value = '  -1.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.25'
$ws.Range("E23").Value = '  -2.35%  '
$ws.Range("E24").Value = '  -3.82%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.90'
$ws.Range("E25").Value = '  -0.92%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.53'
$ws.Range("E26").Value = '  -1.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.14'
$ws.Range("E27").Value = '  +0.94%  '
$ws.Range("E28").Value = '  -3.02%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.80'
$ws.Range("E30").Value = '  -2.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0518'
$ws.Range("E31").Value = '  -1.53%  '
$ws.Range("E32").Value = '  -1.67%  '
$ws.Range("E33").Value = '  +0.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.92'
$ws.Range("E34").Value = '  +6.31%  '
$ws.Range("D35").Value = '1.462.89'
$ws.Range("E35").Value = '  -4.40%  '
$ws.Range("E36").Value = '  -1.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0191'
$ws.Range("E37").Value = '  +0.50%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.635'
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("B39").Value = 'Aave'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '83.37'
$ws.Range("E39").Value = '  -0.51%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.85'
$ws.Range("E40").Value = '  +2.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.32'
$ws.Range("E41").Value = '  -1.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.902'
$ws.Range("E42").Value = '  -1.09%  '
$ws.Range("E43").Value = '  -1.32%  '
$ws.Range("E44").Value = '  -2.46%  '
$ws.Range("E45").Value = '  +1.31%  '
$ws.Range("D46").Value = '1.955.30'
$ws.Range("E46").Value = '  -0.69%  '
$ws.Range("E47").Value = '  -3.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '11.96'
$ws.Range("E48").Value = '  -1.60%  '
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '100.70'
$ws.Range("E50").Value = '  +1.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '49.87'
$ws.Range("E51").Value = '  -4.22%  '
